$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-4 (row 5 gets deleted)
# Columns: A Phone Number, B Original Message, C Timestamp, D Meal Tag,
#          E Items Extracted, F Total Calories, G Total Protein, H Source

$ws.Cells.Item(2, 2).Value = "Had 3 idlis for breakfast"
$ws.Cells.Item(2, 3).Value = "2026-01-14T17:17:48.419079"
$ws.Cells.Item(2, 4).Value = "Evening Snack"
$ws.Cells.Item(2, 5).Value = "3x idli"
$ws.Cells.Item(2, 6).Value = 117
$ws.Cells.Item(2, 7).Value = 6

$ws.Cells.Item(3, 2).Value = "Ate chicken biryani"
$ws.Cells.Item(3, 3).Value = "2026-01-14T17:17:48.412534"
$ws.Cells.Item(3, 4).Value = "Evening Snack"
$ws.Cells.Item(3, 5).Value = "1x biryani"
$ws.Cells.Item(3, 6).Value = 280
$ws.Cells.Item(3, 7).Value = 12

$ws.Cells.Item(4, 2).Value = "I had 2 rotis and dal"
$ws.Cells.Item(4, 3).Value = "2026-01-14T17:17:48.404366"
$ws.Cells.Item(4, 4).Value = "Evening Snack"
$ws.Cells.Item(4, 5).Value = "2x roti, 1x dal"
$ws.Cells.Item(4, 6).Value = 246
$ws.Cells.Item(4, 7).Value = 13.8

# Delete the now-obsolete last row (row 5), shifting cells up
$ws.Range("A5:H5").Delete()
